$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IPC_Dist")
$ws.Activate()
Write-Host $ws.Name
